$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7092821002006531
$ws.Range("B1").Value = 0.6446579694747925
$ws.Range("C1").Value = 0.4804243147373199
$ws.Range("D1").Value = 0.4517580270767212
$ws.Range("E1").Value = 0.4746173620223999
